$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Indirizzo" values between row 2 and row 3 (A2 <-> A3)
$a2 = $ws.Range("A2").Value()
$a3 = $ws.Range("A3").Value()
$ws.Range("A2").Value = $a3
$ws.Range("A3").Value = $a2

# Match the selection left behind in the saved workbook (A4 active cell)
$null = $ws.Range("A4").Select()
